$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44616
$ws.Range("K2").Value = 19000
$ws.Range("L2").Value = 20000
$ws.Range("M2").Value = 19500
$ws.Range("P2").Value = 1500

# Row 3
$ws.Range("D3").Value = 44229
$ws.Range("K3").Value = 44000
$ws.Range("L3").Value = 45000
$ws.Range("M3").Value = 44500
$ws.Range("P3").Value = 3423

# Row 4
$ws.Range("D4").Value = 44406
$ws.Range("J4").Value = 160
$ws.Range("K4").Value = 17000
$ws.Range("L4").Value = 18000
$ws.Range("M4").Value = 17500
$ws.Range("P4").Value = 1346

# Row 5
$ws.Range("D5").Value = 44580
$ws.Range("J5").Value = 160
$ws.Range("K5").Value = 11000
$ws.Range("L5").Value = 12000
$ws.Range("M5").Value = 11500
$ws.Range("P5").Value = 885

# Row 7
$ws.Range("D7").Value = 44592
$ws.Range("K7").Value = 12000
$ws.Range("L7").Value = 13000
$ws.Range("M7").Value = 12500
$ws.Range("P7").Value = 962

# Row 8
$ws.Range("D8").Value = 44397
$ws.Range("J8").Value = 140
$ws.Range("K8").Value = 12500
$ws.Range("L8").Value = 13000
$ws.Range("M8").Value = 12750
$ws.Range("P8").Value = 981

# Row 9
$ws.Range("D9").Value = 44469
$ws.Range("K9").Value = 13000
$ws.Range("L9").Value = 14000
$ws.Range("M9").Value = 13500
$ws.Range("P9").Value = 1038

# Row 10
$ws.Range("D10").Value = 44389
$ws.Range("J10").Value = 120
$ws.Range("K10").Value = 12000
$ws.Range("L10").Value = 13000
$ws.Range("M10").Value = 12500
$ws.Range("P10").Value = 962

# Row 11
$ws.Range("D11").Value = 44320
$ws.Range("K11").Value = 19000
$ws.Range("L11").Value = 20000
$ws.Range("M11").Value = 19500
$ws.Range("P11").Value = 1500

# Row 12
$ws.Range("D12").Value = 44379
$ws.Range("M12").Value = 12667
$ws.Range("P12").Value = 974
